{"js": "// Remove the empty \"separator\" paragraphs (pPr -> spacing w:before=\"40\" only,\n// no runs) that immediately follow a table. These were leftover spacer\n// paragraphs between a table and the next \"\u2022 Archivo: ...\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  items[i].load(\"text,spaceBefore,tableNestingLevel\");\n}\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const isEmpty = (para.text || \"\").trim().length === 0;\n  const isBodyLevel = para.tableNestingLevel === 0;\n  const hasTargetSpacing = para.spaceBefore === 2; // w:before=\"40\" (twentieths) == 2pt\n  const prevIsInsideTable = i > 0 && items[i - 1].tableNestingLevel > 0;\n\n  if (isEmpty && isBodyLevel && hasTargetSpacing && prevIsInsideTable) {\n    toDelete.push(para);\n  }\n}\n\n// Delete in reverse order so earlier indices stay valid while deleting.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the empty \"separator\" paragraphs (pPr -> spacing w:before=\"40\" only,\n# no runs) that immediately follow a table. These were leftover spacer\n# paragraphs between a table and the next \"\u2022 Archivo: ...\" paragraph.\n$d = $word.ActiveDocument\n\n$targets = @()\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n    # A bare paragraph mark (\"\\r\", length 1) with no runs, and\n    # spacing-before of 2pt (w:before=\"40\" twentieths-of-a-point).\n    if ($text -eq [char]13 -and $p.Format.SpaceBefore -eq 2) {\n        $targets += $i\n    }\n}\n\n# Delete in reverse order so earlier indices stay valid while deleting.\nfor ($j = $targets.Count - 1; $j -ge 0; $j--) {\n    $idx = $targets[$j]\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
